# QB Website updated 12/13
# Update the "Occurrence" column (column E) on Sheet1 to reflect the
# newly-added question occurrence number (85) and related list tweaks,
# then leave the selection where the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value  = "19, 21, 23, 25, 27, 31, 33, 37, 39, 41, 45, 49, 53, 57, 61, 65, 69, 73, 81, "
$ws.Range("E6").Value  = "19, 21, 23, 25, 27, 31, 33, 37, 39, 41, 45, 49, 53, 57, 61, 65, 73, 81, 85"
$ws.Range("E8").Value  = "39, 41, 45, 49, 53, 57, 61, 65, 69, 73, 81,"
$ws.Range("E9").Value  = "39, 41, 45, 49, 53, 57, 61, 65, 69, 73, 81, 85"

# Restore the view's active cell/selection to match the saved workbook state.
$ws.Range("G8").Select()
